$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.571000000000001
$ws.Range("D6").Value = -8.445
$ws.Range("D7").Value = -7.703
$ws.Range("D8").Value = -7.876
$ws.Range("D16").Value = -7.922000000000001
$ws.Range("D20").Value = -7.81
$ws.Range("D21").Value = -7.726999999999999
